# 15th September, Data Driven
# Re-order / refresh the login test-data table on Sheet1:
#   - header C1: "result" -> "res"
#   - row 2 becomes the "admin" (valid) case, with A2 turned into a mailto
#     hyperlink (picks up the workbook's themed Hyperlink style automatically)
#   - row 5 becomes a new "none" / 123 (invalid) case with no special styling

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value2 = "res"

# Row 2: admin@admin.com / admin123 / Valid  (was the old row 5 content)
$ws.Range("A2").Value2 = "admin@admin.com"
$ws.Range("B2").Value2 = "admin123"
$ws.Range("C2").Value2 = "Valid"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@admin.com")

# Row 3 stays the same (abc@gmail.com / 1234 / Invalid)
$ws.Range("B3").Value2 = 1234

# Row 4 stays the same (1223.com / asdf / Invalid) - hyperlink untouched

# Row 5: none / 123 / Invalid, with default (unstyled) formatting on A5/B5
$ws.Range("A5").Value2 = "none"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value2 = 123
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value2 = "Invalid"
